$wb = $excel.ActiveWorkbook

# --- TDIL-EN sheet: insert a new row 2 for "Magical Something" ---
$wsTdil = $wb.Worksheets.Item("TDIL-EN")
$wsTdil.Rows.Item(2).Insert()
$wsTdil.Cells.Item(2, 1).Value = "Magical Something"
$wsTdil.Cells.Item(2, 2).Value = 100909000
$wsTdil.Cells.Item(2, 2).ClearFormats()
$wsTdil.Cells.Item(2, 3).Value = ":"
$wsTdil.Cells.Item(2, 5).Value = ","

# --- SPDS-JP sheet: insert "Abyss Actor - Pretty Heroine" label at A22 ---
$wsSpds = $wb.Worksheets.Item("SPDS-JP")
$wsSpds.Cells.Item(22, 1).Value = "Abyss Actor - Pretty Heroine"

# --- Leave TDIL-EN's own selection where it should be, then switch to SPDS-JP last ---
$wsTdil.Activate()
$wsTdil.Cells.Item(5, 6).Select()

# --- Make SPDS-JP the active sheet/tab (must be the LAST activation/selection) ---
$wsSpds.Activate()
$wsSpds.Cells.Item(20, 6).Select()
$excel.ActiveWindow.ScrollRow = 11

# --- TDIL-EN also gained an explicit (portrait) page setup ---
$wsTdil.PageSetup.Orientation = 1
